$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "28.112.70"
$c.Style = $origStyle
$ws.Range("E2").Value = "  -1.00%  "
$c = $ws.Range("D3")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.794.64"
$c.Style = $origStyle
$ws.Range("E3").Value = "  +0.10%  "
$c = $ws.Range("D4")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = $origStyle
$ws.Range("E4").Value = "  -0.18%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "317.19"
$c.Style = $origStyle
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  -0.18%  "
$c = $ws.Range("D7")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.5415"
$c.Style = $origStyle
$ws.Range("E7").Value = "  -0.15%  "
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.3783"
$c.Style = $origStyle
$ws.Range("E8").Value = "  -1.39%  "
$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07451"
$c.Style = $origStyle
$ws.Range("E9").Value = "  -1.60%  "
$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "41.71"
$c.Style = $origStyle
$ws.Range("E10").Value = "  -1.82%  "
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.094"
$c.Style = $origStyle
$ws.Range("E11").Value = "  -2.39%  "
$c = $ws.Range("D12")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = $origStyle
$ws.Range("E12").Value = "  -0.17%  "
$c = $ws.Range("D13")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "20.54"
$c.Style = $origStyle
$ws.Range("E13").Value = "  -2.53%  "
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.115"
$c.Style = $origStyle
$ws.Range("E14").Value = "  -1.07%  "
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.249"
$c.Style = $origStyle
$ws.Range("E15").Value = "  -1.97%  "
$c = $ws.Range("D16")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.782.88"
$c.Style = $origStyle
$ws.Range("E16").Value = "  -0.64%  "
$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "89.13"
$c.Style = $origStyle
$ws.Range("E17").Value = "  -2.81%  "
$c = $ws.Range("D18")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.00001059"
$c.Style = $origStyle
$ws.Range("E18").Value = "  -0.92%  "
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06485"
$c.Style = $origStyle
$ws.Range("E19").Value = "  +0.46%  "
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = $origStyle
$ws.Range("E20").Value = "  -0.20%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "17.27"
$c.Style = $origStyle
$ws.Range("E21").Value = "  -0.26%  "
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.906"
$c.Style = $origStyle
$ws.Range("E22").Value = "  -0.88%  "
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "28.130.23"
$c.Style = $origStyle
$ws.Range("E23").Value = "  -0.96%  "
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "11.17"
$c.Style = $origStyle
$ws.Range("E24").Value = "  -1.44%  "
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.087"
$c.Style = $origStyle
$ws.Range("E25").Value = "  -1.65%  "
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "154.98"
$c.Style = $origStyle
$ws.Range("E26").Value = "  -2.85%  "
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "20.28"
$c.Style = $origStyle
$ws.Range("E27").Value = "  -1.88%  "
$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.994.51"
$c.Style = $origStyle
$ws.Range("E28").Value = "  -0.38%  "
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.289"
$c.Style = $origStyle
$ws.Range("E29").Value = "  -4.47%  "
$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "121.15"
$c.Style = $origStyle
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("E31").Value = "  +0.55%  "
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1056"
$c.Style = $origStyle
$ws.Range("E32").Value = "  +3.58%  "
$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.653"
$c.Style = $origStyle
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("E34").Value = "  -3.06%  "
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.2260"
$c.Style = $origStyle
$ws.Range("E35").Value = "  -2.79%  "
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06507"
$c.Style = $origStyle
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  -1.01%  "
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.019"
$c.Style = $origStyle
$ws.Range("E38").Value = "  -2.34%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.456"
$c.Style = $origStyle
$ws.Range("E39").Value = "  -3.67%  "
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.452"
$c.Style = $origStyle
$ws.Range("E40").Value = "  +4.38%  "
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.6183"
$c.Style = $origStyle
$ws.Range("E41").Value = "  -3.31%  "
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "11.09"
$c.Style = $origStyle
$ws.Range("E42").Value = "  -4.42%  "
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.172"
$c.Style = $origStyle
$ws.Range("E43").Value = "  +1.40%  "
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = $origStyle
$ws.Range("E44").Value = "  -0.24%  "
$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "13.37"
$c.Style = $origStyle
$ws.Range("E45").Value = "  -0.98%  "
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.676"
$c.Style = $origStyle
$ws.Range("E46").Value = "  +0.06%  "
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.5784"
$c.Style = $origStyle
$ws.Range("E47").Value = "  -3.01%  "
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "124.70"
$c.Style = $origStyle
$ws.Range("E48").Value = "  -0.96%  "
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.193"
$c.Style = $origStyle
$ws.Range("E49").Value = "  +3.95%  "
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.924"
$c.Style = $origStyle
$ws.Range("E50").Value = "  -2.89%  "
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.06811"
$c.Style = $origStyle
$ws.Range("E51").Value = "  -1.16%  "
